$d = $word.ActiveDocument

# Each entry is the list of run-texts that make up one new paragraph.
$newParas = @(
    @("09 – Road turn ", "bottom-right"),
    @("10 – ", "Road turn left-bottom"),
    @("11 – Road turn top-left"),
    @("12 – Road turn right-top"),
    @("13 – House top-look"),
    @("14 – House left-look"),
    @("15 – House right-look"),
    @("16 – Road city crossroad", " horizontal"),
    @("17 – Road city crossroad vertical")
)

# Anchor on the last paragraph currently in the document
# ("08 – Road crossroad") and append the new paragraphs after it,
# preserving the same paragraph/run formatting (en-GB language).
$anchor = $d.Paragraphs($d.Paragraphs.Count).Range

foreach ($runTexts in $newParas) {
    $anchor.InsertParagraphAfter()
    $anchor = $d.Paragraphs($d.Paragraphs.Count).Range

    foreach ($runText in $runTexts) {
        $insertAt = $anchor
        $insertAt.Collapse(0)
        $insertAt.InsertAfter($runText)
    }
}
